$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new columns (E, F) mirroring the existing C/D "setProperty"/"saveProperties"
#     test-case pair, by first copying formatting from the analogous existing cells,
#     then writing the new values on top. This keeps the same cell-style (xf) indices
#     that are already used elsewhere on the sheet instead of creating new duplicate styles.

# Row 1 header style (same as D1) -> E1, F1
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# Row 2 style (same as B2, the "json" cell) -> E2
$ws.Range("B2").Copy()
$ws.Range("E2").PasteSpecial(-4122)

# Row 2 style (same as A2, an empty cell) -> F2
$ws.Range("A2").Copy()
$ws.Range("F2").PasteSpecial(-4122)

# Row 3 style (same as B3) -> E3, F3
$ws.Range("B3").Copy()
$ws.Range("E3:F3").PasteSpecial(-4122)

# --- Update existing values that changed
$ws.Range("D2").Value = "src/test/resources/saveproperties/result1.json"
$ws.Range("F3").Value = '{"target":"src/test/resources/saveproperties/result2.json"}'
$ws.Range("C3").Value = '{"result1":"OK"}'

# --- Fill in the new cells' values
$ws.Range("E1").Value = "setProperty"
$ws.Range("F1").Value = "saveProperties"
$ws.Range("E2").Value = "json"
$ws.Range("E3").Value = '{"result2":"NG"}'

# --- Row 3 grows taller (wrapped long JSON text in the new, narrower columns)
$ws.Rows.Item(3).RowHeight = 31.5

# --- Column widths (best fit / resized for the new layout)
$ws.Columns.Item(2).ColumnWidth = 42.428571428571431
$ws.Columns.Item(3).ColumnWidth = 11
$ws.Columns.Item(4).ColumnWidth = 36.142857142857146
$ws.Columns.Item(5).ColumnWidth = 11
$ws.Columns.Item(6).ColumnWidth = 31.285714285714285

# --- Selection moved
$ws.Range("F10").Select()
